$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.139.82"
$ws.Range("E2").Value = "  +1.70%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.139.26"
$ws.Range("E3").Value = "  +2.10%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.80"
$ws.Range("E5").Value = "  +2.88%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.29"
$ws.Range("E6").Value = "  +2.85%  "

# Row 8 - XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  +10.76%  "

# Row 9 - Toncoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +0.10%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.98%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +6.21%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.84%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.679.06"
$ws.Range("E13").Value = "  +2.13%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.03"

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +5.30%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.238.34"
$ws.Range("E16").Value = "  +1.86%  "

# Row 17 - Polkadot -> WrappedEther
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.154.23"
$ws.Range("E17").Value = "  +2.53%  "

# Row 18 - WrappedEther -> Polkadot
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.25"
$ws.Range("E18").Value = "  +6.65%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +4.42%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.19"
$ws.Range("E20").Value = "  +4.85%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.66"
$ws.Range("E21").Value = "  +7.06%  "

# Row 22 - Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23 - LEO
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  -0.41%  "

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.44"
$ws.Range("E24").Value = "  +2.18%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +3.72%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +1.08%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.42%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.00"
$ws.Range("E28").Value = "  +11.17%  "

# Row 29 - PEPE
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0886"
$ws.Range("E29").Value = "  +2.49%  "

# Row 30 - PancakeSwap -> RenderToken
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.17"
$ws.Range("E30").Value = "  +5.86%  "

# Row 31 - RenderToken -> PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.89"
$ws.Range("E31").Value = "  +1.81%  "

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.78"
$ws.Range("E32").Value = "  +4.32%  "

# Row 33 - NEARProtocol
$ws.Range("E33").Value = "  +7.03%  "

# Row 34 - Fetch.AI
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.17"
$ws.Range("E34").Value = "  +4.20%  "

# Row 35 - Monero
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.05"
$ws.Range("E35").Value = "  +2.03%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +4.63%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +11.79%  "

# Row 38 - EnergySwap
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.52"
$ws.Range("E38").Value = "  +0.07%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +7.21%  "

# Row 40 - Maker
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.644.07"
$ws.Range("E40").Value = "  +10.11%  "

# Row 41 - Hedera
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0680"
$ws.Range("E41").Value = "  +3.90%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +4.87%  "

# Row 43 - OKB
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.62"
$ws.Range("E43").Value = "  +5.47%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  +1.33%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +4.82%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.00%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  +11.65%  "

# Row 48 - Cosmos
$ws.Range("E48").Value = "  +4.61%  "

# Row 49 - ONDO
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.979"
$ws.Range("E49").Value = "  +3.81%  "

# Row 50 - InjectiveProtocol
$ws.Range("E50").Value = "  +3.86%  "

# Row 51 - SuiNetwork
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.750"
$ws.Range("E51").Value = "  +0.01%  "

